$d = $word.ActiveDocument

# Insert " Hà Nội" right after "Xây dựng" in the "Kính gửi" addressee line.
$d.Content.Find.Execute("Xây dựng", $true, $false, $false, $false, $false, $true, 1, $false, "Xây dựng Hà Nội", 2)
